# Insert a new data row at row 12 (pushing existing rows 12-23 down to 13-24)
# and populate it with a new weekly price record for Cebollín.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Insert()

$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = 44791
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 100112037
$ws.Cells.Item(12, 7).Value = "Cebollín"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 120
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 8500
$ws.Cells.Item(12, 13).Value = 8250
$ws.Cells.Item(12, 14).Value = "$/docena de atados"
$ws.Cells.Item(12, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(12, 16).Value = 2750
$ws.Cells.Item(12, 17).Value = 3
$ws.Cells.Item(12, 18).Value = "Hortaliza"
